$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 6 data (Array di Linked List / merge di k linked list) ---
$ws.Range("A6").Value = "Array di Linked List"
$ws.Range("B6").Value = "merge di k linked list orinate"
$ws.Range("C6").Value = "creo una mappa con key=valore nodo value=numero occorrenze, iterando per tutte le linke list. Una volta creata la mappa, la ordino mediante tuple e liste e cro k nodi relativi ad una entry quanto è il valore dell'occorrenza."
$ws.Range("D6").Value = "O(n^2)"
$ws.Range("E6").Value = "difficile"

# --- Row heights: rows 2-9 get the taller 28.2pt custom height ---
$ws.Rows.Item(2).RowHeight = 28.2
$ws.Rows.Item(3).RowHeight = 28.2
$ws.Rows.Item(4).RowHeight = 28.2
$ws.Rows.Item(5).RowHeight = 28.2
$ws.Rows.Item(6).RowHeight = 28.2
$ws.Rows.Item(7).RowHeight = 28.2
$ws.Rows.Item(8).RowHeight = 28.2
$ws.Rows.Item(9).RowHeight = 28.2

# --- Highlight column D (Complessità) for rows 2-6 with a new fill + border ---
$colD = $ws.Range("D2:D6")
$colD.Borders.ColorIndex = 1
$colD.Borders.LineStyle = 1
$colD.Interior.ThemeColor = 8
$colD.Interior.TintAndShade = 0.79998168889431442

# --- Selection moves to D11 ---
$ws.Range("D11").Select()
